# The deck's slide master currently carries the "Integral" (Red Violet)
# theme colours. The target revision swaps the slide-master theme back to
# the default "Office Theme" colour palette (dk1/lt1/dk2/lt2/accent1-6/
# hlink/folHlink), which is what ppt/theme/theme2.xml (the theme used by
# the slide master / all slides) must end up containing.
#
# PowerPoint's ThemeColorScheme exposes exactly 12 colour slots, in this
# fixed order:
#   1 dk1  2 lt1  3 dk2  4 lt2  5 accent1 6 accent2 7 accent3 8 accent4
#   9 accent5 10 accent6 11 hlink 12 folHlink
# RGBColor.RGB uses the classic OLE BGR-packed integer (0x00BBGGRR), so a
# hex RRGGBB string is written with the byte order reversed.

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

function Set-ThemeColor($scheme, $index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    $scheme.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

Set-ThemeColor $tcs 1  "000000"   # dk1
Set-ThemeColor $tcs 2  "FFFFFF"   # lt1
Set-ThemeColor $tcs 3  "44546A"   # dk2
Set-ThemeColor $tcs 4  "E7E6E6"   # lt2
Set-ThemeColor $tcs 5  "5B9BD5"   # accent1
Set-ThemeColor $tcs 6  "ED7D31"   # accent2
Set-ThemeColor $tcs 7  "A5A5A5"   # accent3
Set-ThemeColor $tcs 8  "FFC000"   # accent4
Set-ThemeColor $tcs 9  "4472C4"   # accent5
Set-ThemeColor $tcs 10 "70AD47"   # accent6
Set-ThemeColor $tcs 11 "0563C1"   # hlink
Set-ThemeColor $tcs 12 "954F72"   # folHlink
